# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail) right before
#    the "总计" (total) summary sheet.
# 2. Add a new summary row for 2022-Q1 at the top of the "总计" sheet's
#    data table, shifting the existing quarter rows down by one and
#    renumbering the helper index column (A).

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# 1. New "2022-Q1" sheet with the per-fund holding detail
# ------------------------------------------------------------------
# NOTE: worksheet handles here are positional, not stable identities —
# once the new sheet is spliced in "before" it, the variable that used
# to point at "总计" silently starts pointing at the new sheet instead
# (same slot index). So: grab the "before" target, add relative to it,
# rename the new sheet, and ONLY THEN re-fetch "总计" by name (below)
# for the second part of the edit.
$beforeTarget = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeTarget)
$newSheet.Name = "2022-Q1"

# Reuse the existing header / index-column formatting instead of
# re-building fonts/borders by hand, so the new sheet visually matches
# its siblings.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'582003"
$newSheet.Range("C2").Value = "东吴配置优化灵活配置混合"
$newSheet.Range("D2").Value = "'1.04"
$newSheet.Range("E2").Value = "'90.74"
$newSheet.Range("F2").Value = "'5.61"
$newSheet.Range("G2").Value = "'0.0583"
$newSheet.Range("H2").Value = 2

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'002681"
$newSheet.Range("C3").Value = "金鹰元和灵活配置混合A"
$newSheet.Range("D3").Value = "'0.56"
$newSheet.Range("E3").Value = "'81.63"
$newSheet.Range("F3").Value = "'5.00"
$newSheet.Range("G3").Value = "'0.0280"
$newSheet.Range("H3").Value = 4

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'002682"
$newSheet.Range("C4").Value = "金鹰元和灵活配置混合C"
$newSheet.Range("D4").Value = "'0.25"
$newSheet.Range("E4").Value = "'81.63"
$newSheet.Range("F4").Value = "'5.00"
$newSheet.Range("G4").Value = "'0.0125"
$newSheet.Range("H4").Value = 4

# The leading apostrophes force text storage (otherwise "582003" etc
# would be reinterpreted as numbers and lose leading zeros / precision)
# but they also stamp a "quote prefix" style on the cell. Strip that
# back off so the cells fall back to the plain default style, matching
# their neighbours.
$newSheet.Range("B2:B4").ClearFormats()
$newSheet.Range("D2:G4").ClearFormats()

# ------------------------------------------------------------------
# 2. Add the 2022-Q1 row to the "总计" summary sheet
# ------------------------------------------------------------------
# Re-fetch by name now that the sheet collection has shifted (see note
# above) so this grabs the real "总计" sheet, not the new one.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert(-4121)   # xlShiftDown

# Insert() drags odd formatting along with the shifted row; reset the
# new row then restore just the index-column style from its neighbour.
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.1

# Renumber the helper index column for the rows that shifted down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
